$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2753413333333334
$ws.Range("H2").Value = 0.8260240000000001
$ws.Range("I2").Value = 0.006630378892106956
$ws.Range("J2").Value = 0.006630378892106955
$ws.Range("M2").Value = 86.89540866666668
$ws.Range("N2").Value = 260.686226
$ws.Range("O2").Value = 0.319779657009892
$ws.Range("P2").Value = 0.3197796570098919
$ws.Range("Q2").Value = 23.9258976828249
$ws.Range("R2").Value = 215.3330791454241
$ws.Range("S2").Value = 0.00212026028796359
$ws.Range("T2").Value = 0.002120260287963589

$ws.Range("G3").Value = 0.2753413333333334
$ws.Range("H3").Value = 0.8260240000000001
$ws.Range("I3").Value = 0.006630378892106956
$ws.Range("J3").Value = 0.006630378892106955
$ws.Range("O3").Value = 0.1999969065479545
$ws.Range("P3").Value = 0.1999969065479545
$ws.Range("Q3").Value = 14.96375838191556
$ws.Range("R3").Value = 134.67382543724
$ws.Range("S3").Value = 0.001326055267662245
$ws.Range("T3").Value = 0.001326055267662245

$ws.Range("G4").Value = 0.2753413333333334
$ws.Range("H4").Value = 0.8260240000000001
$ws.Range("I4").Value = 0.006630378892106956
$ws.Range("J4").Value = 0.006630378892106955
$ws.Range("M4").Value = 60.92601633333334
$ws.Range("N4").Value = 182.778049
$ws.Range("O4").Value = 0.224210932487692
$ws.Range("P4").Value = 0.224210932487692
$ws.Range("Q4").Value = 16.77545057190845
$ws.Range("R4").Value = 150.979055147176
$ws.Range("S4").Value = 0.001486603434146011
$ws.Range("T4").Value = 0.00148660343414601

$ws.Range("G5").Value = 0.2753413333333334
$ws.Range("H5").Value = 0.8260240000000001
$ws.Range("I5").Value = 0.006630378892106956
$ws.Range("J5").Value = 0.006630378892106955
$ws.Range("M5").Value = 7.809668333333332
$ws.Range("N5").Value = 23.429005
$ws.Range("O5").Value = 0.02873998867505581
$ws.Range("P5").Value = 0.02873998867505581
$ws.Range("Q5").Value = 2.150324491791111
$ws.Range("R5").Value = 19.35292042612
$ws.Range("S5").Value = 0.000190557014270483
$ws.Range("T5").Value = 0.000190557014270483

$ws.Range("G6").Value = 0.2753413333333334
$ws.Range("H6").Value = 0.8260240000000001
$ws.Range("I6").Value = 0.006630378892106956
$ws.Range("J6").Value = 0.006630378892106955
$ws.Range("M6").Value = 61.75795633333333
$ws.Range("N6").Value = 185.273869
$ws.Range("O6").Value = 0.2272725152794058
$ws.Range("P6").Value = 0.2272725152794058
$ws.Range("Q6").Value = 17.00451804076178
$ws.Range("R6").Value = 153.040662366856
$ws.Range("S6").Value = 0.001506902888064628
$ws.Range("T6").Value = 0.001506902888064628

$ws.Range("I7").Value = 0.03952244389885164
$ws.Range("J7").Value = 0.03952244389885164
$ws.Range("M7").Value = 86.89540866666668
$ws.Range("N7").Value = 260.686226
$ws.Range("O7").Value = 0.319779657009892
$ws.Range("P7").Value = 0.3197796570098919
$ws.Range("Q7").Value = 142.617784637436
$ws.Range("R7").Value = 1283.560061736924
$ws.Range("S7").Value = 0.01263847355416748
$ws.Range("T7").Value = 0.01263847355416747

$ws.Range("I8").Value = 0.03952244389885164
$ws.Range("J8").Value = 0.03952244389885164
$ws.Range("O8").Value = 0.1999969065479545
$ws.Range("P8").Value = 0.1999969065479545
$ws.Range("S8").Value = 0.007904366518985408
$ws.Range("T8").Value = 0.007904366518985408

$ws.Range("I9").Value = 0.03952244389885164
$ws.Range("J9").Value = 0.03952244389885164
$ws.Range("M9").Value = 60.92601633333334
$ws.Range("N9").Value = 182.778049
$ws.Range("O9").Value = 0.224210932487692
$ws.Range("P9").Value = 0.224210932487692
$ws.Range("Q9").Value = 99.995311715214
$ws.Range("R9").Value = 899.957805436926
$ws.Range("S9").Value = 0.00886136400075402
$ws.Range("T9").Value = 0.008861364000754018

$ws.Range("I10").Value = 0.03952244389885164
$ws.Range("J10").Value = 0.03952244389885164
$ws.Range("M10").Value = 7.809668333333332
$ws.Range("N10").Value = 23.429005
$ws.Range("O10").Value = 0.02873998867505581
$ws.Range("P10").Value = 0.02873998867505581
$ws.Range("Q10").Value = 12.81768062943
$ws.Range("R10").Value = 115.35912566487
$ws.Range("S10").Value = 0.001135874590063525
$ws.Range("T10").Value = 0.001135874590063525

$ws.Range("I11").Value = 0.03952244389885164
$ws.Range("J11").Value = 0.03952244389885164
$ws.Range("M11").Value = 61.75795633333333
$ws.Range("N11").Value = 185.273869
$ws.Range("O11").Value = 0.2272725152794058
$ws.Range("P11").Value = 0.2272725152794058
$ws.Range("Q11").Value = 101.360739895734
$ws.Range("R11").Value = 912.2466590616059
$ws.Range("S11").Value = 0.008982365234881219
$ws.Range("T11").Value = 0.008982365234881219

$ws.Range("G12").Value = 23.78768866666667
$ws.Range("H12").Value = 71.363066
$ws.Range("I12").Value = 0.5728213302306416
$ws.Range("J12").Value = 0.5728213302306416
$ws.Range("M12").Value = 86.89540866666668
$ws.Range("N12").Value = 260.686226
$ws.Range("O12").Value = 0.319779657009892
$ws.Range("P12").Value = 0.3197796570098919
$ws.Range("Q12").Value = 2067.040927925435
$ws.Range("R12").Value = 18603.36835132892
$ws.Range("S12").Value = 0.1831766085091046
$ws.Range("T12").Value = 0.1831766085091046

$ws.Range("G13").Value = 23.78768866666667
$ws.Range("H13").Value = 71.363066
$ws.Range("I13").Value = 0.5728213302306416
$ws.Range("J13").Value = 0.5728213302306416
$ws.Range("O13").Value = 0.1999969065479545
$ws.Range("P13").Value = 0.1999969065479545
$ws.Range("Q13").Value = 1292.770763339435
$ws.Range("R13").Value = 11634.93687005491
$ws.Range("S13").Value = 0.1145624940508126
$ws.Range("T13").Value = 0.1145624940508126

$ws.Range("G14").Value = 23.78768866666667
$ws.Range("H14").Value = 71.363066
$ws.Range("I14").Value = 0.5728213302306416
$ws.Range("J14").Value = 0.5728213302306416
$ws.Range("M14").Value = 60.92601633333334
$ws.Range("N14").Value = 182.778049
$ws.Range("O14").Value = 0.224210932487692
$ws.Range("P14").Value = 0.224210932487692
$ws.Range("Q14").Value = 1449.289108237582
$ws.Range("R14").Value = 13043.60197413824
$ws.Range("S14").Value = 0.1284328045998523
$ws.Range("T14").Value = 0.1284328045998523

$ws.Range("G15").Value = 23.78768866666667
$ws.Range("H15").Value = 71.363066
$ws.Range("I15").Value = 0.5728213302306416
$ws.Range("J15").Value = 0.5728213302306416
$ws.Range("M15").Value = 7.809668333333332
$ws.Range("N15").Value = 23.429005
$ws.Range("O15").Value = 0.02873998867505581
$ws.Range("P15").Value = 0.02873998867505581
$ws.Range("Q15").Value = 185.7739589032589
$ws.Range("R15").Value = 1671.96563012933
$ws.Range("S15").Value = 0.01646287854365905
$ws.Range("T15").Value = 0.01646287854365905

$ws.Range("G16").Value = 23.78768866666667
$ws.Range("H16").Value = 71.363066
$ws.Range("I16").Value = 0.5728213302306416
$ws.Range("J16").Value = 0.5728213302306416
$ws.Range("M16").Value = 61.75795633333333
$ws.Range("N16").Value = 185.273869
$ws.Range("O16").Value = 0.2272725152794058
$ws.Range("P16").Value = 0.2272725152794058
$ws.Range("Q16").Value = 1469.079037946928
$ws.Range("R16").Value = 13221.71134152235
$ws.Range("S16").Value = 0.1301865445272131
$ws.Range("T16").Value = 0.1301865445272131

$ws.Range("G17").Value = 0.5982033333333333
$ws.Range("H17").Value = 1.79461
$ws.Range("I17").Value = 0.01440508298011203
$ws.Range("J17").Value = 0.01440508298011203
$ws.Range("M17").Value = 86.89540866666668
$ws.Range("N17").Value = 260.686226
$ws.Range("O17").Value = 0.319779657009892
$ws.Range("P17").Value = 0.3197796570098919
$ws.Range("Q17").Value = 51.98112311576223
$ws.Range("R17").Value = 467.8301080418601
$ws.Range("S17").Value = 0.004606452494579258
$ws.Range("T17").Value = 0.004606452494579257

$ws.Range("G18").Value = 0.5982033333333333
$ws.Range("H18").Value = 1.79461
$ws.Range("I18").Value = 0.01440508298011203
$ws.Range("J18").Value = 0.01440508298011203
$ws.Range("O18").Value = 0.1999969065479545
$ws.Range("P18").Value = 0.1999969065479545
$ws.Range("Q18").Value = 32.51008497303889
$ws.Range("R18").Value = 292.59076475735
$ws.Range("S18").Value = 0.002880972034588997
$ws.Range("T18").Value = 0.002880972034588996

$ws.Range("G19").Value = 0.5982033333333333
$ws.Range("H19").Value = 1.79461
$ws.Range("I19").Value = 0.01440508298011203
$ws.Range("J19").Value = 0.01440508298011203
$ws.Range("M19").Value = 60.92601633333334
$ws.Range("N19").Value = 182.778049
$ws.Range("O19").Value = 0.224210932487692
$ws.Range("P19").Value = 0.224210932487692
$ws.Range("Q19").Value = 36.44614605732111
$ws.Range("R19").Value = 328.01531451589
$ws.Range("S19").Value = 0.0032297770875335
$ws.Range("T19").Value = 0.0032297770875335

$ws.Range("G20").Value = 0.5982033333333333
$ws.Range("H20").Value = 1.79461
$ws.Range("I20").Value = 0.01440508298011203
$ws.Range("J20").Value = 0.01440508298011203
$ws.Range("M20").Value = 7.809668333333332
$ws.Range("N20").Value = 23.429005
$ws.Range("O20").Value = 0.02873998867505581
$ws.Range("P20").Value = 0.02873998867505581
$ws.Range("Q20").Value = 4.671769629227777
$ws.Range("R20").Value = 42.04592666305
$ws.Range("S20").Value = 0.0004140019217116591
$ws.Range("T20").Value = 0.0004140019217116591

$ws.Range("G21").Value = 0.5982033333333333
$ws.Range("H21").Value = 1.79461
$ws.Range("I21").Value = 0.01440508298011203
$ws.Range("J21").Value = 0.01440508298011203
$ws.Range("M21").Value = 61.75795633333333
$ws.Range("N21").Value = 185.273869
$ws.Range("O21").Value = 0.2272725152794058
$ws.Range("P21").Value = 0.2272725152794058
$ws.Range("Q21").Value = 36.94381533845444
$ws.Range("R21").Value = 332.49433804609
$ws.Range("S21").Value = 0.003273879441698621
$ws.Range("T21").Value = 0.003273879441698621

$ws.Range("G22").Value = 15.22474833333333
$ws.Range("H22").Value = 45.674245
$ws.Range("I22").Value = 0.3666207639982877
$ws.Range("J22").Value = 0.3666207639982877
$ws.Range("M22").Value = 86.89540866666668
$ws.Range("N22").Value = 260.686226
$ws.Range("O22").Value = 0.319779657009892
$ws.Range("P22").Value = 0.3197796570098919
$ws.Range("Q22").Value = 1322.960728272152
$ws.Range("R22").Value = 11906.64655444937
$ws.Range("S22").Value = 0.117237862164077
$ws.Range("T22").Value = 0.117237862164077

$ws.Range("G23").Value = 15.22474833333333
$ws.Range("H23").Value = 45.674245
$ws.Range("I23").Value = 0.3666207639982877
$ws.Range("J23").Value = 0.3666207639982877
$ws.Range("O23").Value = 0.1999969065479545
$ws.Range("P23").Value = 0.1999969065479545
$ws.Range("Q23").Value = 827.4073954950638
$ws.Range("R23").Value = 7446.666559455575
$ws.Range("S23").Value = 0.07332301867590525
$ws.Range("T23").Value = 0.07332301867590522

$ws.Range("G24").Value = 15.22474833333333
$ws.Range("H24").Value = 45.674245
$ws.Range("I24").Value = 0.3666207639982877
$ws.Range("J24").Value = 0.3666207639982877
$ws.Range("M24").Value = 60.92601633333334
$ws.Range("N24").Value = 182.778049
$ws.Range("O24").Value = 0.224210932487692
$ws.Range("P24").Value = 0.224210932487692
$ws.Range("Q24").Value = 927.5832656275561
$ws.Range("R24").Value = 8348.249390648005
$ws.Range("S24").Value = 0.08220038336540615
$ws.Range("T24").Value = 0.08220038336540614

$ws.Range("G25").Value = 15.22474833333333
$ws.Range("H25").Value = 45.674245
$ws.Range("I25").Value = 0.3666207639982877
$ws.Range("J25").Value = 0.3666207639982877
$ws.Range("M25").Value = 7.809668333333332
$ws.Range("N25").Value = 23.429005
$ws.Range("O25").Value = 0.02873998867505581
$ws.Range("P25").Value = 0.02873998867505581
$ws.Range("Q25").Value = 118.9002349418028
$ws.Range("R25").Value = 1070.102114476225
$ws.Range("S25").Value = 0.0105366766053511
$ws.Range("T25").Value = 0.0105366766053511

$ws.Range("G26").Value = 15.22474833333333
$ws.Range("H26").Value = 45.674245
$ws.Range("I26").Value = 0.3666207639982877
$ws.Range("J26").Value = 0.3666207639982877
$ws.Range("M26").Value = 61.75795633333333
$ws.Range("N26").Value = 185.273869
$ws.Range("O26").Value = 0.2272725152794058
$ws.Range("P26").Value = 0.2272725152794058
$ws.Range("Q26").Value = 940.2493427559893
$ws.Range("R26").Value = 8462.244084803904
$ws.Range("S26").Value = 0.0833228231875483
$ws.Range("T26").Value = 0.08332282318754827
